$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 110.91
$ws.Range("F8").Value = 50.92

# Row 9
$ws.Range("C9").Value = 110.44
$ws.Range("F9").Value = 50.38

# Row 10
$ws.Range("C10").Value = 110.25
$ws.Range("F10").Value = 50.11

# Row 11
$ws.Range("C11").Value = 111.27
$ws.Range("F11").Value = 51.07

# Row 12
$ws.Range("C12").Value = 110.51
$ws.Range("F12").Value = 50.24

# Row 13
$ws.Range("C13").Value = 154.9
$ws.Range("F13").Value = 50.66

# Row 14
$ws.Range("C14").Value = 298.11
$ws.Range("F14").Value = 51.1

# Row 16
$ws.Range("C16").Value = 858.55
$ws.Range("F16").Value = 50.14

# Row 17
$ws.Range("C17").Value = 200.43
$ws.Range("F17").Value = 49.68

# Row 18 (was "-" string, now numeric)
$ws.Range("C18").Value = 104.85
$ws.Range("F18").Value = 49.44

# Row 19
$ws.Range("C19").Value = 105.09
$ws.Range("F19").Value = 49.46

# Row 20
$ws.Range("C20").Value = 105.33
$ws.Range("F20").Value = 49.47

# Row 21
$ws.Range("C21").Value = 105.58
$ws.Range("F21").Value = 49.49

# Row 22
$ws.Range("C22").Value = 105.82
$ws.Range("F22").Value = 49.5

# Row 23 (C23 was numeric, now becomes "-" string; F23 numeric update)
$ws.Range("C23").Value = "-"
$ws.Range("F23").Value = 49.51

# Row 24
$ws.Range("F24").Value = 49.53

# Row 25
$ws.Range("F25").Value = 49.55

# Row 26
$ws.Range("F26").Value = 49.56

# Row 27
$ws.Range("F27").Value = 49.58

# Row 28
$ws.Range("F28").Value = 49.68

# Row 29
$ws.Range("F29").Value = 50.29
